# Update the thresholds_summary table and the "Chosen thresholds" notes
# to reflect the final-logic re-run of the dedup analysis.

$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# --- Table cell updates (row index includes the header row = row 1) ---

# threshold 0.55 -> row_count/unique_count 7 -> 6
$tbl.Cell(13, 2).Range.Text = "6"
$tbl.Cell(13, 3).Range.Text = "6"

# threshold 0.65 -> row_count/unique_count 7 -> 9
$tbl.Cell(15, 2).Range.Text = "9"
$tbl.Cell(15, 3).Range.Text = "9"

# threshold 0.7 -> row_count/unique_count 7 -> 9
$tbl.Cell(16, 2).Range.Text = "9"
$tbl.Cell(16, 3).Range.Text = "9"

# threshold 0.75 -> row_count/unique_count 7 -> 9
$tbl.Cell(17, 2).Range.Text = "9"
$tbl.Cell(17, 3).Range.Text = "9"

# threshold 0.8 -> row_count/unique_count 7 -> 9
$tbl.Cell(18, 2).Range.Text = "9"
$tbl.Cell(18, 3).Range.Text = "9"

# threshold 0.85 -> row_count/unique_count 7 -> 10
$tbl.Cell(19, 2).Range.Text = "10"
$tbl.Cell(19, 3).Range.Text = "10"

# threshold 0.9 -> row_count/unique_count 7 -> 10
$tbl.Cell(20, 2).Range.Text = "10"
$tbl.Cell(20, 3).Range.Text = "10"

# threshold 0.95 -> row_count 9 -> 10, unique_count 7 -> 10,
#                   repeated_count 4 -> 0, no_repeats_bool False -> True
$tbl.Cell(21, 2).Range.Text = "10"
$tbl.Cell(21, 3).Range.Text = "10"
$tbl.Cell(21, 4).Range.Text = "0"
$tbl.Cell(21, 5).Range.Text = "True"

# threshold 1.0 -> unique_count 7 -> 10, repeated_count 5 -> 0,
#                  no_repeats_bool False -> True (row_count was already 10)
$tbl.Cell(22, 3).Range.Text = "10"
$tbl.Cell(22, 4).Range.Text = "0"
$tbl.Cell(22, 5).Range.Text = "True"

# --- "Chosen thresholds" paragraph text updates ---

$d.Content.Find.Execute(
    "  - Full Analysis threshold = 0.90.", $true, $false, $false, $false,
    $false, $true, 1, $false,
    "  - Full Analysis threshold = 1.00.", 2)

$d.Content.Find.Execute(
    "  - Core-Level Analysis: NONE found (no threshold yields no repeats & ≤9 rows?).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "  - Core-Level Analysis threshold = 0.80.", 2)
